$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for every existing data row
# (rows 2 through 515) from 45177 (2023-09-08) to 45178 (2023-09-09).
for ($r = 2; $r -le 515; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}

# Row 515 gains an explicit row height (matches every other data row).
$ws.Rows.Item(515).RowHeight = 15

# Append the new record as row 516.
$ws.Range("A516").Value = "A 41945-2023"

$ws.Range("B516").Value = 45177
$ws.Range("B516").NumberFormat = "YYYY-MM-DD"

$ws.Range("C516").Value = 45178
$ws.Range("C516").NumberFormat = "YYYY-MM-DD"

$ws.Range("D516").Value = "UPPSALA LÄN"
$ws.Range("E516").Value = "ENKÖPING"

$ws.Range("G516").Value = 0.6
$ws.Range("H516").Value = 0
$ws.Range("I516").Value = 0
$ws.Range("J516").Value = 0
$ws.Range("K516").Value = 0
$ws.Range("L516").Value = 0
$ws.Range("M516").Value = 0
$ws.Range("N516").Value = 0
$ws.Range("O516").Value = 0
$ws.Range("P516").Value = 0
$ws.Range("Q516").Value = 0

$ws.Range("R516").WrapText = $true
